$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: Title slide
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Date placeholder: refresh cached "datetime1" field text (12/29/23 -> 1/8/24 run)
$s1DateShape = $s1.Shapes.Item(3)
$s1DateShape.TextFrame.TextRange.Text = "2024-01-08"

# Title: "Lecture 10: Conditional Computation" -> "Lecture 10: Scaling Law"
# (split across 3 runs, matching how the edit was actually typed)
$s1Title = $s1.Shapes.Item(1)
$s1Title.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/> (drops fontScale)
$titleRange = $s1Title.TextFrame.TextRange
$titleRange.Text = "Lecture 10"
[void]$titleRange.InsertAfter(": Scaling ")
[void]$titleRange.InsertAfter("Law")

# ---------------------------------------------------------------------------
# Slide 2: "What you need to cover"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Date placeholder: refresh cached "datetime1" field text
$s2DateShape = $s2.Shapes.Item(3)
$s2DateShape.TextFrame.TextRange.Text = "2024-01-08"

# Content placeholder: replace bullet content with the new "scaling law" notes
$content = $s2.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange

$lines = @(
    "In machine learning, a neural scaling law is a scaling law relating parameters of a family of neural networks.",
    "The parameters are interested in are size of the model, size of the training data, cost of training and test loss.",
    "Different laws",
    "OpenAI",
    "Chinchilla (Hoffmann et al. 2022)",
    "Emergent Ability (Wei et al. 2022)",
    "Debate about Emergent Ability",
    "Other papers."
)
$contentRange.Text = [string]::Join("`r", $lines)

# Sub-bullets (level 2) start at paragraph 4
for ($i = 4; $i -le 8; $i++) {
    $para = $contentRange.Paragraphs($i, 1)
    $para.IndentLevel = 2
}

# "OpenAI" + " GPT (Kaplan et al. 2020)" as two runs in the same paragraph
$openAiPara = $contentRange.Paragraphs(4, 1)
[void]$openAiPara.InsertAfter(" GPT (Kaplan et al. 2020)")
